# Applies the commit's content updates to the FHIR StructureDefinition export.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
# Date regenerated
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
# FHIR Version downgraded from 4.3.0 (R4B) to 4.0.1 (R4)
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet --------------------------------------------------------
# Extension root element: drop the "unless an empty Parameters resource ..."
# carve-out from the ele-1 invariant text (now matches the plain ele-1 wording
# used elsewhere in the sheet).
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id: type changed from "id" to "string"
$wsElem.Range("K3").Value = "string" + [char]10

# Extension.value[x] definition: FHIR spec link moves from R4B to R4
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
